$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compartments")
Write-Host $ws.Name
